$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts old A->B, old B->C),
# picking up the neighboring column's per-cell formatting automatically.
$ws.Columns("A:A").Insert()

# New "Car_PR_ID" column header + ids for the existing 4 rows.
$ws.Range("A1").Value = "Car_PR_ID"
$ws.Range("A2").Value = "Car_PR_01"
$ws.Range("A3").Value = "Car_PR_02"
$ws.Range("A4").Value = "Car_PR_03"

# Copy A4's cell format down into the new rows (5-11) without touching
# their (still empty) values, so they pick up the shared "s=8" style.
$ws.Range("A4").Copy()
$ws.Range("A5:A11").PasteSpecial(-4122)

# Fill in the new PR id rows + the extra reviewer-comment rows.
$ws.Range("A5").Value = "Car_PR_04"
$ws.Range("C5").Value = "there is no related SRS ID's column "

$ws.Range("A6").Value = "Car_PR_05"
$ws.Range("C6").Value = "there is no column for the reviewer "

$ws.Range("A7").Value = "Car_PR_06"
$ws.Range("C7").Value = "please try to work on the last version of the project cause there is  new updates"

$ws.Range("A8").Value = "Car_PR_07"
$ws.Range("A9").Value = "Car_PR_08"
$ws.Range("A10").Value = "Car_PR_09"
$ws.Range("A11").Value = "Car_PR_10"

# Match the row height used on the rest of the sheet.
$ws.Range("A5:A11").RowHeight = 21

# Match the saved selection from the source workbook.
$ws.Range("C10").Select()
